$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.90240433333333
$ws.Range("H2").Value = 83.707213
$ws.Range("I2").Value = 0.2174736967445081
$ws.Range("J2").Value = 0.2174736967445081
$ws.Range("M2").Value = 8.658142333333332
$ws.Range("N2").Value = 25.974427
$ws.Range("O2").Value = 0.1344393815966713
$ws.Range("P2").Value = 0.1344393815966713
$ws.Range("Q2").Value = 241.5829881602168
$ws.Range("R2").Value = 2174.246893441951
$ws.Range("S2").Value = 0.02923702930387368
$ws.Range("T2").Value = 0.02923702930387368

$ws.Range("G3").Value = 27.90240433333333
$ws.Range("H3").Value = 83.707213
$ws.Range("I3").Value = 0.2174736967445081
$ws.Range("J3").Value = 0.2174736967445081
$ws.Range("O3").Value = 0.04342195384682981
$ws.Range("P3").Value = 0.04342195384682981
$ws.Range("Q3").Value = 78.02777160596445
$ws.Range("R3").Value = 702.24994445368
$ws.Range("S3").Value = 0.009443132822939493
$ws.Range("T3").Value = 0.009443132822939493

$ws.Range("G4").Value = 27.90240433333333
$ws.Range("H4").Value = 83.707213
$ws.Range("I4").Value = 0.2174736967445081
$ws.Range("J4").Value = 0.2174736967445081
$ws.Range("M4").Value = 4.152187333333333
$ws.Range("N4").Value = 12.456562
$ws.Range("O4").Value = 0.06447312551305152
$ws.Range("P4").Value = 0.06447312551305154
$ws.Range("Q4").Value = 115.8560098424118
$ws.Range("R4").Value = 1042.704088581706
$ws.Range("S4").Value = 0.01402120894599597
$ws.Range("T4").Value = 0.01402120894599598

$ws.Range("G5").Value = 27.90240433333333
$ws.Range("H5").Value = 83.707213
$ws.Range("I5").Value = 0.2174736967445081
$ws.Range("J5").Value = 0.2174736967445081
$ws.Range("M5").Value = 48.795048
$ws.Range("N5").Value = 146.385144
$ws.Range("O5").Value = 0.7576655390434474
$ws.Range("P5").Value = 0.7576655390434474
$ws.Range("Q5").Value = 1361.499158760408
$ws.Range("R5").Value = 12253.49242884367
$ws.Range("S5").Value = 0.1647723256716989
$ws.Range("T5").Value = 0.1647723256716989

$ws.Range("G6").Value = 64.92210766666668
$ws.Range("I6").Value = 0.5060083921817455
$ws.Range("J6").Value = 0.5060083921817455
$ws.Range("M6").Value = 8.658142333333332
$ws.Range("N6").Value = 25.974427
$ws.Range("O6").Value = 0.1344393815966713
$ws.Range("P6").Value = 0.1344393815966713
$ws.Range("Q6").Value = 562.1048487579912
$ws.Range("R6").Value = 5058.943638821921
$ws.Range("S6").Value = 0.06802745532763978
$ws.Range("T6").Value = 0.06802745532763978

$ws.Range("G7").Value = 64.92210766666668
$ws.Range("I7").Value = 0.5060083921817455
$ws.Range("J7").Value = 0.5060083921817455
$ws.Range("O7").Value = 0.04342195384682981
$ws.Range("P7").Value = 0.04342195384682981
$ws.Range("Q7").Value = 181.5516443914756
$ws.Range("S7").Value = 0.02197187305142432
$ws.Range("T7").Value = 0.02197187305142432

$ws.Range("G8").Value = 64.92210766666668
$ws.Range("I8").Value = 0.5060083921817455
$ws.Range("J8").Value = 0.5060083921817455
$ws.Range("M8").Value = 4.152187333333333
$ws.Range("N8").Value = 12.456562
$ws.Range("O8").Value = 0.06447312551305152
$ws.Range("P8").Value = 0.06447312551305154
$ws.Range("Q8").Value = 269.5687531068363
$ws.Range("R8").Value = 2426.118777961526
$ws.Range("S8").Value = 0.03262394257979108
$ws.Range("T8").Value = 0.03262394257979109

$ws.Range("G9").Value = 64.92210766666668
$ws.Range("I9").Value = 0.5060083921817455
$ws.Range("J9").Value = 0.5060083921817455
$ws.Range("M9").Value = 48.795048
$ws.Range("N9").Value = 146.385144
$ws.Range("O9").Value = 0.7576655390434474
$ws.Range("P9").Value = 0.7576655390434474
$ws.Range("Q9").Value = 3167.877359856168
$ws.Range("R9").Value = 28510.89623870551
$ws.Range("S9").Value = 0.3833851212228904
$ws.Range("T9").Value = 0.3833851212228904

$ws.Range("G10").Value = 19.423329
$ws.Range("H10").Value = 58.269987
$ws.Range("I10").Value = 0.1513870672309258
$ws.Range("J10").Value = 0.1513870672309258
$ws.Range("M10").Value = 8.658142333333332
$ws.Range("N10").Value = 25.974427
$ws.Range("O10").Value = 0.1344393815966713
$ws.Range("P10").Value = 0.1344393815966713
$ws.Range("Q10").Value = 168.169947069161
$ws.Range("R10").Value = 1513.529523622449
$ws.Range("S10").Value = 0.02035238370025936
$ws.Range("T10").Value = 0.02035238370025936

$ws.Range("G11").Value = 19.423329
$ws.Range("H11").Value = 58.269987
$ws.Range("I11").Value = 0.1513870672309258
$ws.Range("J11").Value = 0.1513870672309258
$ws.Range("O11").Value = 0.04342195384682981
$ws.Range("P11").Value = 0.04342195384682981
$ws.Range("Q11").Value = 54.31643312648
$ws.Range("R11").Value = 488.84789813832
$ws.Range("S11").Value = 0.006573522246308182
$ws.Range("T11").Value = 0.006573522246308183

$ws.Range("G12").Value = 19.423329
$ws.Range("H12").Value = 58.269987
$ws.Range("I12").Value = 0.1513870672309258
$ws.Range("J12").Value = 0.1513870672309258
$ws.Range("M12").Value = 4.152187333333333
$ws.Range("N12").Value = 12.456562
$ws.Range("O12").Value = 0.06447312551305152
$ws.Range("P12").Value = 0.06447312551305154
$ws.Range("Q12").Value = 80.64930064496599
$ws.Range("R12").Value = 725.843705804694
$ws.Range("S12").Value = 0.009760397386632248
$ws.Range("T12").Value = 0.009760397386632253

$ws.Range("G13").Value = 19.423329
$ws.Range("H13").Value = 58.269987
$ws.Range("I13").Value = 0.1513870672309258
$ws.Range("J13").Value = 0.1513870672309258
$ws.Range("M13").Value = 48.795048
$ws.Range("N13").Value = 146.385144
$ws.Range("O13").Value = 0.7576655390434474
$ws.Range("P13").Value = 0.7576655390434474
$ws.Range("Q13").Value = 947.762270874792
$ws.Range("R13").Value = 8529.860437873127
$ws.Range("S13").Value = 0.114700763897726
$ws.Range("T13").Value = 0.114700763897726

$ws.Range("G14").Value = 16.05459166666667
$ws.Range("H14").Value = 48.163775
$ws.Range("I14").Value = 0.1251308438428206
$ws.Range("J14").Value = 0.1251308438428206
$ws.Range("M14").Value = 8.658142333333332
$ws.Range("N14").Value = 25.974427
$ws.Range("O14").Value = 0.1344393815966713
$ws.Range("P14").Value = 0.1344393815966713
$ws.Range("Q14").Value = 139.0029397535472
$ws.Range("R14").Value = 1251.026457781925
$ws.Range("S14").Value = 0.01682251326489844
$ws.Range("T14").Value = 0.01682251326489844

$ws.Range("G15").Value = 16.05459166666667
$ws.Range("H15").Value = 48.163775
$ws.Range("I15").Value = 0.1251308438428206
$ws.Range("J15").Value = 0.1251308438428206
$ws.Range("O15").Value = 0.04342195384682981
$ws.Range("P15").Value = 0.04342195384682981
$ws.Range("Q15").Value = 44.89591638155556
$ws.Range("R15").Value = 404.063247434
$ws.Range("S15").Value = 0.005433425726157823
$ws.Range("T15").Value = 0.005433425726157823

$ws.Range("G16").Value = 16.05459166666667
$ws.Range("H16").Value = 48.163775
$ws.Range("I16").Value = 0.1251308438428206
$ws.Range("J16").Value = 0.1251308438428206
$ws.Range("M16").Value = 4.152187333333333
$ws.Range("N16").Value = 12.456562
$ws.Range("O16").Value = 0.06447312551305152
$ws.Range("P16").Value = 0.06447312551305154
$ws.Range("Q16").Value = 66.66167216017222
$ws.Range("R16").Value = 599.95504944155
$ws.Range("S16").Value = 0.008067576600632221
$ws.Range("T16").Value = 0.008067576600632223

$ws.Range("G17").Value = 16.05459166666667
$ws.Range("H17").Value = 48.163775
$ws.Range("I17").Value = 0.1251308438428206
$ws.Range("J17").Value = 0.1251308438428206
$ws.Range("M17").Value = 48.795048
$ws.Range("N17").Value = 146.385144
$ws.Range("O17").Value = 0.7576655390434474
$ws.Range("P17").Value = 0.7576655390434474
$ws.Range("Q17").Value = 783.3845709954001
$ws.Range("R17").Value = 7050.4611389586
$ws.Range("S17").Value = 0.09480732825113207
$ws.Range("T17").Value = 0.09480732825113207
